$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8633374571800232
$ws.Range("B1").Value = 1.43573522567749
$ws.Range("C1").Value = 4.354164123535156
$ws.Range("D1").Value = 1.346722841262817
$ws.Range("E1").Value = 0.7535381317138672
